$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$c = $t.Cell(1,1)
$c.Range.Text = "52 x 43" + $nl + "  4    3" + $nl + "  ----" + $nl + "5|    |" + $nl + "2|    |"
$c = $t.Cell(1,2)
$c.Range.Text = "24 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"
$c = $t.Cell(1,3)
$c.Range.Text = "65 x 20" + $nl + "  2    0" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |"
$c = $t.Cell(2,1)
$c.Range.Text = "58 x 94" + $nl + "  9    4" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"
$c = $t.Cell(2,2)
$c.Range.Text = "37 x 32" + $nl + "  3    2" + $nl + "  ----" + $nl + "3|    |" + $nl + "7|    |"
$c = $t.Cell(2,3)
$c.Range.Text = "48 x 92" + $nl + "  9    2" + $nl + "  ----" + $nl + "4|    |" + $nl + "8|    |"
$c = $t.Cell(3,1)
$c.Range.Text = "83 x 20" + $nl + "  2    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "3|    |"
$c = $t.Cell(3,2)
$c.Range.Text = "64 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "6|    |" + $nl + "4|    |"
$c = $t.Cell(3,3)
$c.Range.Text = "26 x 53" + $nl + "  5    3" + $nl + "  ----" + $nl + "2|    |" + $nl + "6|    |"
$c = $t.Cell(4,1)
$c.Range.Text = "83 x 20" + $nl + "  2    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "3|    |"
$c = $t.Cell(4,2)
$c.Range.Text = "69 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$c = $t.Cell(4,3)
$c.Range.Text = "54 x 99" + $nl + "  9    9" + $nl + "  ----" + $nl + "5|    |" + $nl + "4|    |"
$c = $t.Cell(5,1)
$c.Range.Text = "68 x 60" + $nl + "  6    0" + $nl + "  ----" + $nl + "6|    |" + $nl + "8|    |"
$c = $t.Cell(5,2)
$c.Range.Text = "50 x 53" + $nl + "  5    3" + $nl + "  ----" + $nl + "5|    |" + $nl + "0|    |"
$c = $t.Cell(5,3)
$c.Range.Text = "43 x 42" + $nl + "  4    2" + $nl + "  ----" + $nl + "4|    |" + $nl + "3|    |"
